$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sciences")

$ws.Range("N54").Value = "https://ar.wikipedia.org/wiki/أشعرية"

$ws.Range("N7").Value = "https://ar.m.wikisource.org/wiki/الإبانة_عن_أصول_الديانة"
$ws.Hyperlinks.Add($ws.Range("N7"), "https://ar.m.wikisource.org/wiki/الإبانة_عن_أصول_الديانة")
$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H59").Value = "طبعة @ 2023/06/05 م - 1444/11/16 هـ"

$ws.PageSetup.Zoom = 55
